$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header updates (columns B:E)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates: B2 cleared, C2:E2 new values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 4.3508091582938189
$ws.Range("D2").Value = 2.3066508007200173
$ws.Range("E2").Value = 3.817409235023514

# Row 3 updates
$ws.Range("B3").Value = 1.2988960876363769
$ws.Range("C3").Value = 7.900352088866569
$ws.Range("D3").Value = 6.319657363627865
$ws.Range("E3").Value = 10.676319606241041

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
